# Update "想去人数" (F column) figures on both the "展览" sheet and the
# "全部类型" sheet, which carries a duplicate copy of the same exhibition
# rows (with one extra row inserted, shifting some row numbers by +1
# starting at row 23).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Worksheets.Item(1)) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F5").Value = 1136
$ws1.Range("F6").Value = 14217
$ws1.Range("F7").Value = 16044
$ws1.Range("F8").Value = 13
$ws1.Range("F9").Value = 70
$ws1.Range("F11").Value = 197
$ws1.Range("F24").Value = 6399
$ws1.Range("F26").Value = 1109
$ws1.Range("F27").Value = 1
$ws1.Range("F28").Value = 5653
$ws1.Range("F29").Value = 87
$ws1.Range("F31").Value = 150
$ws1.Range("F32").Value = 4662

# --- Sheet "全部类型" (Worksheets.Item(4)) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F5").Value = 1136
$ws4.Range("F6").Value = 14217
$ws4.Range("F7").Value = 16044
$ws4.Range("F8").Value = 13
$ws4.Range("F9").Value = 70
$ws4.Range("F11").Value = 197
$ws4.Range("F25").Value = 6399
$ws4.Range("F27").Value = 1109
$ws4.Range("F28").Value = 1
$ws4.Range("F30").Value = 5653
$ws4.Range("F31").Value = 87
$ws4.Range("F33").Value = 150
$ws4.Range("F34").Value = 4662
